# Update the "dSF" column (F) values for the specified rows.
# Mapping of row number -> new value for column F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = 3
    5  = 0
    7  = -3
    10 = 0
    14 = -1
    15 = 1
    16 = 3
    21 = 0
    22 = -2
    29 = 2
    30 = 4
    33 = -2
    35 = 3
    36 = -1
    37 = 2
    39 = -1
    44 = 0
    48 = -1
    49 = 2
    51 = 3
    55 = -1
    57 = 1
    58 = -1
    62 = 6
    64 = 0
    65 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
